# Updates the "想去人数" (number of people interested) column F values
# across all four sheets to reflect the newer scraped counts
# (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1249
$ws.Range("F4").Value = 54
$ws.Range("F5").Value = 5531
$ws.Range("F6").Value = 1768
$ws.Range("F7").Value = 6329
$ws.Range("F8").Value = 137
$ws.Range("F9").Value = 1904
$ws.Range("F15").Value = 48
$ws.Range("F16").Value = 7844
$ws.Range("F17").Value = 7844
$ws.Range("F22").Value = 1739
$ws.Range("F28").Value = 171
$ws.Range("F29").Value = 1717
$ws.Range("F30").Value = 795
$ws.Range("F31").Value = 363
$ws.Range("F34").Value = 74
$ws.Range("F35").Value = 87
$ws.Range("F36").Value = 3912

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 5
$ws.Range("F14").Value = 26

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9538
$ws.Range("F3").Value = 2270
$ws.Range("F5").Value = 262

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9538
$ws.Range("F3").Value = 2270
$ws.Range("F5").Value = 1249
$ws.Range("F7").Value = 54
$ws.Range("F10").Value = 5531
$ws.Range("F11").Value = 262
$ws.Range("F12").Value = 1768
$ws.Range("F13").Value = 6329
$ws.Range("F14").Value = 137
$ws.Range("F15").Value = 1904
$ws.Range("F22").Value = 48
$ws.Range("F23").Value = 7844
$ws.Range("F24").Value = 7844
$ws.Range("F29").Value = 1739
$ws.Range("F34").Value = 171
$ws.Range("F35").Value = 1717
$ws.Range("F36").Value = 795
$ws.Range("F37").Value = 5
$ws.Range("F38").Value = 363
$ws.Range("F40").Value = 26
$ws.Range("F45").Value = 3912
